# Update numeric "F" column values (e.g. view/sales counts) on the
# "展览" sheet and their mirrored rows on the "全部类型" sheet.
#
# 展览 (sheet1):  F7: 591->592, F12: 3058->3059, F20: 71->72,
#                 F23: 434->435, F25: 4692->4694
# 全部类型 (sheet4, mirrors the same source rows): F15: 591->592,
#                 F21: 3058->3059, F29: 71->72, F34: 434->435,
#                 F36: 4692->4694

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F7").Value = 592
$wsExhibit.Range("F12").Value = 3059
$wsExhibit.Range("F20").Value = 72
$wsExhibit.Range("F23").Value = 435
$wsExhibit.Range("F25").Value = 4694

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F15").Value = 592
$wsAll.Range("F21").Value = 3059
$wsAll.Range("F29").Value = 72
$wsAll.Range("F34").Value = 435
$wsAll.Range("F36").Value = 4694
